$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns F, G, H with same style as existing headers (copy from E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header style (bold, border, centered) from an existing header cell to the new ones
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Boolean values for rows 2-12, columns F (KNN_Outliers_MAD), G (SVM_Outliers_MAD), H (RF_Outliers_MAD)
$values = @{
    2  = @($false, $false, $false)
    3  = @($false, $false, $false)
    4  = @($false, $false, $false)
    5  = @($false, $false, $false)
    6  = @($false, $false, $false)
    7  = @($false, $false, $false)
    8  = @($false, $true,  $true)
    9  = @($false, $true,  $true)
    10 = @($false, $false, $false)
    11 = @($false, $false, $false)
    12 = @($false, $false, $false)
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 6).Value = $v[0]
    $ws.Cells.Item($row, 7).Value = $v[1]
    $ws.Cells.Item($row, 8).Value = $v[2]
}
